$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
Write-Host "Name: $($ws.Name)"
Write-Host "Before insert dim: $($ws.UsedRange.Address())"

$ws.Rows("1:2").Insert()

Write-Host "After insert dim: $($ws.UsedRange.Address())"
Write-Host ($ws.Range("A3").Value)
Write-Host ($ws.Range("A4").Value)
Write-Host ($ws.Range("B3").Comment)
